$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. The former column B ("QueueName" /
# "yt_queue") shifts to column C, and the new column B inherits the
# formatting of column A.
$ws.Columns("B:B").Insert()

# Update headers
$ws.Range("A1").Value = "TargetChannelName"
$ws.Range("B1").Value = "TargetChannelID"

# Update data row
$ws.Range("A2").Value = "CharlieKirk"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "UCfaIu2jO-fppCQV_lchCRIQ"

# Resize columns to match the final layout
$ws.Columns("A").ColumnWidth = 23.166666666666668
$ws.Columns("B").ColumnWidth = 31
$ws.Columns("C").ColumnWidth = 17.666666666666668

# Move the active selection
$ws.Range("C5").Select()
